# Update p-values in the "ecological" column of the first table
# (Cod: p-values for trend lines in Fig 8) on the Status rows:
#   GDP 2016:      0.02 -> 0.07
#   OHI economic:  0.87 -> 0.92
#   OHI fisheries: 0.11 -> 0.88
#   Readiness:     0.82 -> 1.00
#   Vulnerability: 0.52 -> 0.96
# TAC and RECOVERY tables / columns are untouched.

$d = $word.ActiveDocument

$table = $d.Tables.Item(1)

$table.Cell(2, 2).Range.Text = "0.07"
$table.Cell(3, 2).Range.Text = "0.92"
$table.Cell(4, 2).Range.Text = "0.88"
$table.Cell(5, 2).Range.Text = "1.00"
$table.Cell(6, 2).Range.Text = "0.96"
